$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Mass_Fractions")
$ws.Range("B2").Value = [double]"7.5943634049407224E-2"
$ws.Range("C2").Value = [double]"4.219131060181059E-3"
$ws.Range("D2").Value = 0.11461536284168035
$ws.Range("E2").Value = [double]"6.9090002709211633E-3"
$ws.Range("F2").Value = [double]"2.1785783678512083E-2"
$ws.Range("H2").Value = 0.287732562868526
$ws.Range("I2").Value = 0.42427863914517178
$ws.Range("K2").Value = [double]"7.1339141443824789E-4"
$ws.Range("M2").Value = [double]"6.3802494671161961E-2"
$ws.Range("B3").Value = [double]"8.3004685819327947E-2"
$ws.Range("C3").Value = [double]"4.613479262540833E-3"
$ws.Range("D3").Value = 0.12930083699642483
$ws.Range("E3").Value = [double]"7.3423152153845803E-3"
$ws.Range("F3").Value = [double]"1.9997835516329554E-2"
$ws.Range("H3").Value = 0.226179314894383
$ws.Range("I3").Value = 0.45382904297834392
$ws.Range("K3").Value = [double]"8.1355483070470942E-4"
$ws.Range("M3").Value = [double]"7.4918934486560618E-2"
$ws.Range("B4").Value = [double]"9.5229880905198364E-2"
$ws.Range("C4").Value = [double]"3.5492355593062272E-3"
$ws.Range("D4").Value = 0.12321852723611326
$ws.Range("E4").Value = [double]"9.0806051440410207E-3"
$ws.Range("F4").Value = [double]"4.1455826886504184E-2"
$ws.Range("H4").Value = [double]"9.9510439123667999E-2"
$ws.Range("I4").Value = 0.56269098102888626
$ws.Range("K4").Value = [double]"6.0435468852042924E-4"
$ws.Range("M4").Value = [double]"6.4660149427762326E-2"
$ws.Range("B5").Value = [double]"9.6312507739279588E-2"
$ws.Range("C5").Value = [double]"2.5060546242810615E-3"
$ws.Range("D5").Value = 0.11793103702537136
$ws.Range("E5").Value = [double]"9.9689018086355219E-3"
$ws.Range("F5").Value = [double]"6.162055808232024E-2"
$ws.Range("H5").Value = [double]"4.8868980519427528E-2"
$ws.Range("I5").Value = 0.61524452862486489
$ws.Range("K5").Value = [double]"4.0080177420546672E-4"
$ws.Range("M5").Value = [double]"4.7146629801614277E-2"
$ws.Range("B6").Value = [double]"9.7556872438494702E-2"
$ws.Range("C6").Value = [double]"1.3095539678176804E-3"
$ws.Range("D6").Value = [double]"9.9742275916372927E-2"
$ws.Range("E6").Value = [double]"1.0404302189858937E-2"
$ws.Range("F6").Value = [double]"9.6624594538377948E-2"
$ws.Range("H6").Value = [double]"2.5268929829495763E-2"
$ws.Range("I6").Value = 0.6433463342058855
$ws.Range("K6").Value = [double]"2.1644459621760089E-4"
$ws.Range("M6").Value = [double]"2.5530692317478801E-2"
$ws.Range("B7").Value = [double]"8.2299028721604536E-2"
$ws.Range("C7").Value = [double]"8.0540559924903161E-4"
$ws.Range("D7").Value = [double]"8.0770735330052232E-2"
$ws.Range("E7").Value = [double]"1.1181294686439296E-2"
$ws.Range("F7").Value = 0.1171752365211685
$ws.Range("H7").Value = [double]"2.9954045284091529E-3"
$ws.Range("I7").Value = 0.68668036373075658
$ws.Range("K7").Value = [double]"1.3179098418896262E-4"
$ws.Range("M7").Value = [double]"1.796073989813185E-2"
$ws.Range("B8").Value = [double]"5.8413052794829778E-2"
$ws.Range("C8").Value = [double]"1.5870375684867312E-4"
$ws.Range("D8").Value = [double]"4.8034117724348957E-2"
$ws.Range("E8").Value = [double]"1.1796185776060016E-2"
$ws.Range("F8").Value = 0.15422841289381825
$ws.Range("H8").Value = [double]"3.1803172224408201E-4"
$ws.Range("I8").Value = 0.72386786315594076
$ws.Range("K8").Value = [double]"2.3188282481681027E-5"
$ws.Range("M8").Value = [double]"3.1604438934276749E-3"
$ws.Range("B9").Value = [double]"3.9495533522106316E-2"
$ws.Range("D9").Value = [double]"3.269526254054541E-2"
$ws.Range("E9").Value = [double]"1.1855314193272237E-2"
$ws.Range("F9").Value = 0.17995647932240258
$ws.Range("H9").Value = [double]"5.7785961129601892E-4"
$ws.Range("I9").Value = 0.7347829650987584
$ws.Range("K9").Value = [double]"4.4261618007156399E-6"
$ws.Range("M9").Value = [double]"6.321595498183332E-4"
$ws.Range("B10").Value = [double]"1.4459886264583246E-2"
$ws.Range("D10").Value = [double]"1.1865009655018815E-2"
$ws.Range("E10").Value = [double]"1.2145176273342695E-2"
$ws.Range("F10").Value = 0.21306943527723599
$ws.Range("I10").Value = 0.74846049252981917

$ws = $wb.Worksheets.Item("Uncertainties")
$ws.Range("B2").Value = [double]"9.2050117066182456E-3"
$ws.Range("C2").Value = [double]"3.3609724666102563E-4"
$ws.Range("D2").Value = [double]"2.2979121317208685E-2"
$ws.Range("E2").Value = [double]"6.1370432448422418E-4"
$ws.Range("F2").Value = [double]"2.0569928391430204E-3"
$ws.Range("H2").Value = [double]"7.2536173916128294E-2"
$ws.Range("I2").Value = [double]"2.6771951190910557E-2"
$ws.Range("K2").Value = [double]"6.2995427529744831E-5"
$ws.Range("M2").Value = [double]"4.0338204234419686E-3"
$ws.Range("B3").Value = [double]"1.1204041524548663E-2"
$ws.Range("C3").Value = [double]"6.7687925502405832E-4"
$ws.Range("D3").Value = [double]"3.3499305283460984E-2"
$ws.Range("E3").Value = [double]"1.2097500210485378E-3"
$ws.Range("F3").Value = [double]"1.1758540626339927E-2"
$ws.Range("H3").Value = 0.1177652160048267
$ws.Range("I3").Value = [double]"6.3857950699712868E-2"
$ws.Range("K3").Value = [double]"1.3663296247674346E-4"
$ws.Range("M3").Value = [double]"9.9742287911236752E-3"
$ws.Range("B4").Value = [double]"1.0224437265368716E-2"
$ws.Range("C4").Value = [double]"5.095689670170224E-4"
$ws.Range("D4").Value = [double]"3.687671521963614E-2"
$ws.Range("E4").Value = [double]"1.0677610586162034E-3"
$ws.Range("F4").Value = [double]"1.0389077804728837E-2"
$ws.Range("H4").Value = [double]"3.8170275331183827E-2"
$ws.Range("I4").Value = [double]"3.8248577916091514E-2"
$ws.Range("K4").Value = [double]"7.9536493069314116E-5"
$ws.Range("M4").Value = [double]"7.7613996789671173E-3"
$ws.Range("B5").Value = [double]"5.2019544803832947E-3"
$ws.Range("C5").Value = [double]"3.1704651892104342E-4"
$ws.Range("D5").Value = [double]"2.7956938557151732E-2"
$ws.Range("E5").Value = [double]"9.3528708405208271E-4"
$ws.Range("F5").Value = [double]"6.2542256543869131E-3"
$ws.Range("H5").Value = [double]"1.5915136903460986E-2"
$ws.Range("I5").Value = [double]"2.9060917574679154E-2"
$ws.Range("K5").Value = [double]"8.7667049512448435E-5"
$ws.Range("M5").Value = [double]"3.9314019883051167E-3"
$ws.Range("B6").Value = [double]"6.2220254396788419E-3"
$ws.Range("C6").Value = [double]"3.5719583796350293E-4"
$ws.Range("D6").Value = [double]"2.7455374002231081E-2"
$ws.Range("E6").Value = [double]"8.9376183043998487E-4"
$ws.Range("F6").Value = [double]"9.614610501094906E-3"
$ws.Range("H6").Value = [double]"1.4634778950769102E-2"
$ws.Range("I6").Value = [double]"2.5826827578832583E-2"
$ws.Range("K6").Value = [double]"6.6019016451554019E-5"
$ws.Range("M6").Value = [double]"5.8193001581905105E-3"
$ws.Range("B7").Value = [double]"6.2086006262574911E-3"
$ws.Range("C7").Value = [double]"5.7882600176676193E-5"
$ws.Range("D7").Value = [double]"1.9771611650357908E-2"
$ws.Range("E7").Value = [double]"8.3623237068937764E-4"
$ws.Range("F7").Value = [double]"9.3267757297111828E-3"
$ws.Range("H7").Value = [double]"1.0052958459337001E-3"
$ws.Range("I7").Value = [double]"1.8482430287661514E-2"
$ws.Range("K7").Value = [double]"7.5216538753842362E-6"
$ws.Range("M7").Value = [double]"2.1012878523818468E-3"
$ws.Range("B8").Value = [double]"2.0521867080031237E-3"
$ws.Range("C8").Value = [double]"1.2272873485124463E-5"
$ws.Range("D8").Value = [double]"9.1907984360299691E-3"
$ws.Range("E8").Value = [double]"8.6480127968323339E-4"
$ws.Range("F8").Value = [double]"4.0252175709759967E-3"
$ws.Range("H8").Value = [double]"2.6626694117635953E-4"
$ws.Range("I8").Value = [double]"1.7317782005123915E-2"
$ws.Range("K8").Value = [double]"3.0873822898442821E-6"
$ws.Range("M8").Value = [double]"2.2691971644007769E-4"
$ws.Range("B9").Value = [double]"8.9335471994275218E-4"
$ws.Range("D9").Value = [double]"6.1962150112454028E-3"
$ws.Range("E9").Value = [double]"8.3420555875280821E-4"
$ws.Range("F9").Value = [double]"3.5401875121342271E-3"
$ws.Range("H9").Value = [double]"5.9691805208705527E-5"
$ws.Range("I9").Value = [double]"1.4329578028222857E-2"
$ws.Range("K9").Value = [double]"2.3025386894147772E-7"
$ws.Range("M9").Value = [double]"1.3992253335041258E-5"
$ws.Range("B10").Value = [double]"3.1854513254508596E-4"
$ws.Range("D10").Value = [double]"2.2477631024692001E-3"
$ws.Range("E10").Value = [double]"8.5232713784051716E-4"
$ws.Range("F10").Value = [double]"4.0465404584002947E-3"
$ws.Range("I10").Value = [double]"1.4082130486613919E-2"
